# Generate Report for Handback
#
# Adds "Latest Target File" (F) / "Latest Handback File" (G) hyperlink data
# for the two rows on each language sheet (zh-cn, de-de), refreshes the
# "Latest Handback DateTime" (H) for both sheets, and flips the Overview /
# per-language "Status" column text from "Ready for handoff" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Status text: every cell currently showing "Ready for handoff"
#    (Overview!B2:C3, and the Status column C2:C3 on each language sheet)
#    flips to the handed-back message.
# ---------------------------------------------------------------------
$statusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# 2) zh-cn sheet: add "Latest Target File" (F) / "Latest Handback File" (G)
#    hyperlinks for rows 2 and 3, and refresh "Latest Handback DateTime" (H).
# ---------------------------------------------------------------------
$zhMdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/29aa192d37dd6d43052c337bb127863c12e89fc1/e2e/33ee0c4d-bde2-45bb-87c6-bac0e56b5171.md"
$zhMdDisplay = "33ee0c4d-bde2-45bb-87c6-bac0e56b5171.md"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a4651b0b7df41961520df22c9751394a74e7598/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/33ee0c4d-bde2-45bb-87c6-bac0e56b5171.d3764f2e973b14a36cf01934dbfd9ff00b7d3d5e.zh-cn.xlf"
$zhXlfDisplay = "33ee0c4d-bde2-45bb-87c6-bac0e56b5171.d3764f2e973b14a36cf01934dbfd9ff00b7d3d5e.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhMdTarget, [System.Type]::Missing, [System.Type]::Missing, $zhMdDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfTarget, [System.Type]::Missing, [System.Type]::Missing, $zhXlfDisplay)

$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhMdTarget, [System.Type]::Missing, [System.Type]::Missing, $zhMdDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfTarget, [System.Type]::Missing, [System.Type]::Missing, $zhXlfDisplay)

$wsZh.Range("H2").Value = "2016-03-22 05:05:15"
$wsZh.Range("H3").Value = "2016-03-22 05:05:15"

# ---------------------------------------------------------------------
# 3) de-de sheet: same shape of edit, using the de-de handoff package.
# ---------------------------------------------------------------------
$deMdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/29aa192d37dd6d43052c337bb127863c12e89fc1/e2e/33ee0c4d-bde2-45bb-87c6-bac0e56b5171.md"
$deMdDisplay = "33ee0c4d-bde2-45bb-87c6-bac0e56b5171.md"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c962e20c0d2fe1e15582abd3478afae1e40ffa2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/33ee0c4d-bde2-45bb-87c6-bac0e56b5171.d3764f2e973b14a36cf01934dbfd9ff00b7d3d5e.de-de.xlf"
$deXlfDisplay = "33ee0c4d-bde2-45bb-87c6-bac0e56b5171.d3764f2e973b14a36cf01934dbfd9ff00b7d3d5e.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deMdTarget, [System.Type]::Missing, [System.Type]::Missing, $deMdDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfTarget, [System.Type]::Missing, [System.Type]::Missing, $deXlfDisplay)

$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deMdTarget, [System.Type]::Missing, [System.Type]::Missing, $deMdDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfTarget, [System.Type]::Missing, [System.Type]::Missing, $deXlfDisplay)

$wsDe.Range("H2").Value = "2016-03-22 05:05:21"
$wsDe.Range("H3").Value = "2016-03-22 05:05:21"
